# remove type piracy on mean: introduce new colmean to fix it
# update tests (this workbook) to match: formulas/text referencing mean(...)
# are renamed to colmean(...) on the Transformations sheet.

$wb = $excel.ActiveWorkbook

$wsTransformations = $wb.Worksheets.Item("Transformations")
$wsTransformations.Range("B3").Value = "hcat(first_group.flo,second_group.flo).-colmean(third_group.flo)"
$wsTransformations.Range("B4").Value = "hcat(first_group.OD,second_group.OD).-colmean(third_group.OD)"
$wsTransformations.Range("B5").Select()

$wsSamples = $wb.Worksheets.Item("Samples")
$wsSamples.Activate()
